# Issue 50917: test coverage for assay batch files, assay result file old and new root
#
# Adds a new "resultFileField" column (F) to the GenericAssay_Run2 result
# sheet: a bold header in F1 and a sample file-name value ("help.jpg") in
# F2, matching the existing header/body look of the rest of the table.
# Also moves the sheet's active selection off the old B41 cell onto F3
# (where the new data now lives) and clears the stale "scrolled down to
# row 15" view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell - same bold Arial styling as the other header cells
# (A1:E1), picked up automatically since F1 already carried that header
# style in the template.
$ws.Range("F1").Value = "resultFileField"
$ws.Range("F1").Font.Name = "Arial"
$ws.Range("F1").Font.Size = 10
$ws.Range("F1").Font.Bold = $true

# Sample result-file value for the first data row, styled like the plain
# (non-bold) data cells.
$ws.Range("F2").Value = "help.jpg"
$ws.Range("F2").Font.Name = "Arial"
$ws.Range("F2").Font.Size = 10
$ws.Range("F2").Font.Bold = $false

# Refresh the view: select the newly-added cell instead of the old,
# now-irrelevant B41 selection, and reset scroll position.
[void]$ws.Activate()
[void]$ws.Range("A1").Select()
[void]$ws.Range("F3").Select()
